$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Mapping of (row, col) -> new text, following document order.
$updates = @(
    @{ Row = 1;  Col = 1; Text = "40÷7=" },
    @{ Row = 1;  Col = 2; Text = "83÷6=" },
    @{ Row = 1;  Col = 3; Text = "87÷5=" },
    @{ Row = 1;  Col = 4; Text = "21÷2=" },
    @{ Row = 1;  Col = 5; Text = "10÷3=" },

    @{ Row = 5;  Col = 1; Text = "48÷8=" },
    @{ Row = 5;  Col = 2; Text = "55÷8=" },
    @{ Row = 5;  Col = 3; Text = "65÷7=" },
    @{ Row = 5;  Col = 4; Text = "79÷5=" },
    @{ Row = 5;  Col = 5; Text = "86÷5=" },

    @{ Row = 9;  Col = 1; Text = "65÷6=" },
    @{ Row = 9;  Col = 2; Text = "81÷2=" },
    @{ Row = 9;  Col = 3; Text = "19÷5=" },
    @{ Row = 9;  Col = 4; Text = "77÷6=" },
    @{ Row = 9;  Col = 5; Text = "22÷7=" },

    @{ Row = 13; Col = 1; Text = "64÷6=" },
    @{ Row = 13; Col = 2; Text = "63÷2=" },
    @{ Row = 13; Col = 3; Text = "91÷5=" },
    @{ Row = 13; Col = 4; Text = "18÷3=" },
    @{ Row = 13; Col = 5; Text = "88÷2=" },

    @{ Row = 17; Col = 1; Text = "33÷2=" },
    @{ Row = 17; Col = 2; Text = "58÷8=" },
    @{ Row = 17; Col = 3; Text = "78÷4=" },
    @{ Row = 17; Col = 4; Text = "63÷8=" },
    @{ Row = 17; Col = 5; Text = "27÷4=" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $r = $cell.Range
    # Trim the trailing cell-mark/paragraph-mark characters from the range
    # so only the visible text is replaced, then set the new text.
    $r.End = $r.End - 1
    $r.Text = $u.Text
}
